$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.461.21'
$ws.Range('E2').Value = '  -2.71%  '
$ws.Range('D3').Value = '1.774.56'
$ws.Range('E3').Value = '  -1.65%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.84'
$ws.Range('E5').Value = '  -0.90%  '
$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4284'
$ws.Range('E7').Value = '  +1.78%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3630'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07172'
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8437'
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.49'
$ws.Range('E11').Value = '  +1.82%  '
$ws.Range('D12').Value = '1.807.75'
$ws.Range('E12').Value = '  +0.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.254'
$ws.Range('E13').Value = '  -1.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.439'
$ws.Range('E14').Value = '  +1.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06907'
$ws.Range('E15').Value = '  -0.21%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.004'
$ws.Range('E16').Value = '  -0.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '78.89'
$ws.Range('E17').Value = '  -2.53%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008701'
$ws.Range('E18').Value = '  -0.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.93'
$ws.Range('E20').Value = '  -1.05%  '
$ws.Range('D21').Value = '26.470.06'
$ws.Range('E21').Value = '  -2.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.114'
$ws.Range('E22').Value = '  +0.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.08'
$ws.Range('E23').Value = '  +2.11%  '
$ws.Range('D24').Value = '1.995.03'
$ws.Range('E24').Value = '  -1.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.34'
$ws.Range('E25').Value = '  -0.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.871'
$ws.Range('E26').Value = '  -4.59%  '
$ws.Range('E27').Value = '  -0.73%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.068'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '113.88'
$ws.Range('E29').Value = '  +0.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.789'
$ws.Range('E30').Value = '  +4.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08912'
$ws.Range('E31').Value = '  +0.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7262'
$ws.Range('E32').Value = '  -1.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.120'
$ws.Range('E33').Value = '  +1.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.334'
$ws.Range('E34').Value = '  -2.79%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.748'
$ws.Range('E35').Value = '  -7.09%  '
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.002'
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.106'
$ws.Range('E37').Value = '  +3.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05145'
$ws.Range('E38').Value = '  -1.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01892'
$ws.Range('E39').Value = '  -0.41%  '
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.1613'
$ws.Range('E40').Value = '  -1.44%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4922'
$ws.Range('E41').Value = '  -1.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.627'
$ws.Range('E42').Value = '  -4.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.311'
$ws.Range('E43').Value = '  +0.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.004'
$ws.Range('E44').Value = '  -2.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '104.92'
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.003'
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.16'
$ws.Range('E47').Value = '  -1.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.636'
$ws.Range('E48').Value = '  +2.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06200'
$ws.Range('E49').Value = '  -3.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4482'
$ws.Range('E50').Value = '  -1.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.711'
$ws.Range('E51').Value = '  +2.61%  '
